$wb = $excel.ActiveWorkbook

$wsFlights = $wb.Worksheets.Item("Flights")
$wsTestData = $wb.Worksheets.Item("TestData")

# Update the Execute flag for Test_2 from "N" to "Y"
$wsFlights.Range("B3").Value = "Y"

# Switch active sheet to Flights and select C5
$wsFlights.Activate()
$wsFlights.Range("C5").Select()

# Restore TestData's own selection (not active) to G7
$wsTestData.Range("G7").Select()

# Re-activate Flights as the final active sheet
$wsFlights.Activate()
